# Adding MPA test automation upload file
#
# This script applies the following changes to the workbook:
#  1. "Field List" sheet: remove the blank spacer row (old row 2), shifting all
#     field rows up by one, and clear the "Mandatory for sheet" marker that was
#     on the "Document Date in Document" (BLDAT) row.
#  2. Shared text updates: several field labels gain a "(YYYY-MM-DD)" /
#     length-suffix annotation, and two new annotated header labels are used.
#  3. "Data" sheet: header row 5 gets the new annotated labels for the
#     "Revenue from asset sale" / "Amount posted" columns, and a handful of
#     asset / sub-asset numbers used as sample data are bumped to new values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Field List sheet - drop the empty spacer row under the header row.
# ---------------------------------------------------------------------------
$fieldList = $wb.Worksheets.Item("Field List")
$fieldList.Rows.Item(2).Delete()

# After the delete, "Document Date in Document" (BLDAT) is row 4; it is no
# longer a mandatory field, so clear the "Mandatory for sheet" marker.
$fieldList.Range("C4").ClearContents()

# Update field descriptions to call out the expected date format.
$fieldList.Range("B4").Value = "Document Date in Document (YYYY-MM-DD)"
$fieldList.Range("B5").Value = "Posting Date in the Document (YYYY-MM-DD)"
$fieldList.Range("B6").Value = "Asset Value Date (YYYY-MM-DD)"

# ---------------------------------------------------------------------------
# 2) Data sheet - annotate the header row and refresh sample asset numbers.
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# Row 5 holds the human readable / annotated column headers.
$data.Range("C5").Value = "*Document Date in Document (YYYY-MM-DD) (8)"
$data.Range("D5").Value = "*Posting Date in the Document (YYYY-MM-DD) (8)"
$data.Range("E5").Value = "*Asset Value Date (YYYY-MM-DD) (8)"
$data.Range("P5").Value = "Revenue from asset sale (23)"
$data.Range("V5").Value = "Amount posted (23)"

# Sample asset numbers (ANLN1/PANL1) and sub-numbers (ANLN2/PANL2) were
# renumbered for the refreshed upload file.
$mainAssetRows  = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)
foreach ($r in $mainAssetRows) {
    $cell = $data.Range("K$r").Value
    if ($cell -eq 60000391) {
        $data.Range("K$r").Value = 60000409
    }
}

$subAssetRows = @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)
foreach ($r in $subAssetRows) {
    $cell = $data.Range("L$r").Value
    if ($cell -eq 278) {
        $data.Range("L$r").Value = 285
    }
}

$partnerMainRows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $partnerMainRows) {
    $cell = $data.Range("N$r").Value
    if ($cell -eq 60000392) {
        $data.Range("N$r").Value = 60000410
    }
}

$partnerSubRows = @(8, 13, 18, 23, 28)
foreach ($r in $partnerSubRows) {
    $cell = $data.Range("O$r").Value
    if ($cell -eq 279) {
        $data.Range("O$r").Value = 286
    }
}
